# Auto-generated Excel COM-interop script for GEORGIA_2017.xlsx edit
# - Rename header columns to English snake_case names
# - Title-case Spanish connector words (de/del/la/las/el/los/y) in
#   state/municipality names throughout the data rows
# - Fix a handful of 1-ULP floating point percentage values
# - Remove trailing footer/metadata rows (1610-1614)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: translate Spanish labels to English column names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case Spanish connector words in state/municipality names ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B9').Value = 'San Francisco De Los Romo'
$ws.Range('B10').Value = 'San José De Gracia'
$ws.Range('B14').Value = 'Playas De Rosarito'
$ws.Range('B36').Value = 'Amatenango De La Frontera'
$ws.Range('B39').Value = 'Bejucal De Ocampo'
$ws.Range('B41').Value = 'Benemérito De Las Américas'
$ws.Range('B47').Value = 'Chiapa De Corzo'
$ws.Range('B50').Value = 'Comitán De Domínguez'
$ws.Range('B76').Value = 'Mazapa De Madero'
$ws.Range('B81').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B88').Value = 'Salto De Agua'
$ws.Range('B89').Value = 'San Cristóbal De Las Casas'
$ws.Range('B127').Value = 'Guadalupe Y Calvo'
$ws.Range('B129').Value = 'Hidalgo Del Parral'
$ws.Range('B162').Value = 'San Juan De Sabinas'
$ws.Range('B173').Value = 'Villa De Álvarez'
$ws.Range('A175').Value = 'Ciudad De México'
$ws.Range('B179').Value = 'Cuajimalpa De Morelos'
$ws.Range('B193').Value = 'Coneto De Comonfort'
$ws.Range('B207').Value = 'Nombre De Dios'
$ws.Range('B210').Value = 'Pánuco De Coronado'
$ws.Range('B215').Value = 'San Juan De Guadalupe'
$ws.Range('B216').Value = 'San Juan Del Río'
$ws.Range('A225').Value = 'Estado De México'
$ws.Range('B225').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B228').Value = 'Almoloya De Alquisiras'
$ws.Range('B229').Value = 'Almoloya De Juárez'
$ws.Range('B230').Value = 'Almoloya Del Río'
$ws.Range('B236').Value = 'Atizapán De Zaragoza'
$ws.Range('B242').Value = 'Chapa De Mota'
$ws.Range('B245').Value = 'Coacalco De Berriozábal'
$ws.Range('B251').Value = 'Ecatepec De Morelos'
$ws.Range('B258').Value = 'Ixtapan De La Sal'
$ws.Range('B259').Value = 'Ixtapan Del Oro'
$ws.Range('B273').Value = 'Naucalpan De Juárez'
$ws.Range('B282').Value = 'San Felipe Del Progreso'
$ws.Range('B284').Value = 'San Simón De Guerrero'
$ws.Range('B286').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B296').Value = 'Tenango Del Valle'
$ws.Range('B308').Value = 'Tlalnepantla De Baz'
$ws.Range('B313').Value = 'Valle De Bravo'
$ws.Range('B314').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B315').Value = 'Villa De Allende'
$ws.Range('B316').Value = 'Villa Del Carbón'
$ws.Range('B328').Value = 'San Miguel De Allende'
$ws.Range('B329').Value = 'Apaseo El Alto'
$ws.Range('B330').Value = 'Apaseo El Grande'
$ws.Range('B338').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B342').Value = 'Jaral Del Progreso'
$ws.Range('B350').Value = 'Purísima Del Rincón'
$ws.Range('B354').Value = 'San Diego De La Unión'
$ws.Range('B356').Value = 'San Francisco Del Rincón'
$ws.Range('B358').Value = 'San Luis De La Paz'
$ws.Range('B360').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B361').Value = 'Silao De La Victoria'
$ws.Range('B366').Value = 'Valle De Santiago'
$ws.Range('B372').Value = 'Acapulco De Juárez'
$ws.Range('B374').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B375').Value = 'Alcozauca De Guerrero'
$ws.Range('B379').Value = 'Atenango Del Río'
$ws.Range('B380').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B382').Value = 'Atoyac De Álvarez'
$ws.Range('B383').Value = 'Ayutla De Los Libres'
$ws.Range('B386').Value = 'Buenavista De Cuéllar'
$ws.Range('B387').Value = 'Chilapa De Álvarez'
$ws.Range('B388').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B389').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B394').Value = 'Coyuca De Benítez'
$ws.Range('B395').Value = 'Coyuca De Catalán'
$ws.Range('B399').Value = 'Cuetzala Del Progreso'
$ws.Range('B400').Value = 'Cutzamala De Pinzón'
$ws.Range('B406').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B407').Value = 'Iguala De La Independencia'
$ws.Range('B409').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B410').Value = 'Zihuatanejo De Azueta'
$ws.Range('B412').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B415').Value = 'Mártir De Cuilapan'
$ws.Range('B428').Value = 'Taxco De Alarcón'
$ws.Range('B430').Value = 'Técpan De Galeana'
$ws.Range('B432').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B434').Value = 'Tixtla De Guerrero'
$ws.Range('B437').Value = 'Tlapa De Comonfort'
$ws.Range('B448').Value = 'Agua Blanca De Iturbide'
$ws.Range('B454').Value = 'Atotonilco De Tula'
$ws.Range('B455').Value = 'Atotonilco El Grande'
$ws.Range('B461').Value = 'Cuautepec De Hinojosa'
$ws.Range('B466').Value = 'Huasca De Ocampo'
$ws.Range('B470').Value = 'Huejutla De Reyes'
$ws.Range('B473').Value = 'Jacala De Ledezma'
$ws.Range('B480').Value = 'Mineral De La Reforma'
$ws.Range('B481').Value = 'Mineral Del Chico'
$ws.Range('B482').Value = 'Mineral Del Monte'
$ws.Range('B483').Value = 'Mixquiahuala De Juárez'
$ws.Range('B484').Value = 'Molango De Escamilla'
$ws.Range('B486').Value = 'Nopala De Villagrán'
$ws.Range('B487').Value = 'Omitlán De Juárez'
$ws.Range('B488').Value = 'Pachuca De Soto'
$ws.Range('B491').Value = 'Progreso De Obregón'
$ws.Range('B496').Value = 'Santiago De Anaya'
$ws.Range('B497').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B501').Value = 'Tenango De Doria'
$ws.Range('B503').Value = 'Tepehuacán De Guerrero'
$ws.Range('B504').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B505').Value = 'Tezontepec De Aldama'
$ws.Range('B514').Value = 'Tula De Allende'
$ws.Range('B515').Value = 'Tulancingo De Bravo'
$ws.Range('B519').Value = 'Zacualtipán De Ángeles'
$ws.Range('B523').Value = 'Acatlán De Juárez'
$ws.Range('B524').Value = 'Ahualulco De Mercado'
$ws.Range('B529').Value = 'Atemajac De Brizuela'
$ws.Range('B532').Value = 'Atotonilco El Alto'
$ws.Range('B534').Value = 'Autlán De Navarro'
$ws.Range('B548').Value = 'Encarnación De Díaz'
$ws.Range('B553').Value = 'Huejuquilla El Alto'
$ws.Range('B554').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B555').Value = 'Ixtlahuacán Del Río'
$ws.Range('B558').Value = 'Jilotlán De Los Dolores'
$ws.Range('B563').Value = 'La Manzanilla De La Paz'
$ws.Range('B564').Value = 'Lagos De Moreno'
$ws.Range('B569').Value = 'Ojuelos De Jalisco'
$ws.Range('B574').Value = 'San Diego De Alejandría'
$ws.Range('B575').Value = 'San Juan De Los Lagos'
$ws.Range('B578').Value = 'San Miguel El Alto'
$ws.Range('B579').Value = 'Santa María De Los Ángeles'
$ws.Range('B580').Value = 'Santa María Del Oro'
$ws.Range('B583').Value = 'Talpa De Allende'
$ws.Range('B584').Value = 'Tamazula De Gordiano'
$ws.Range('B589').Value = 'Tepatitlán De Morelos'
$ws.Range('B592').Value = 'Tizapán El Alto'
$ws.Range('B593').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B602').Value = 'Unión De San Antonio'
$ws.Range('B603').Value = 'Valle De Guadalupe'
$ws.Range('B604').Value = 'Valle De Juárez'
$ws.Range('B608').Value = 'Yahualica De González Gallo'
$ws.Range('B609').Value = 'Zacoalco De Torres'
$ws.Range('B612').Value = 'Zapotlán Del Rey'
$ws.Range('B613').Value = 'Zapotlán El Grande'
$ws.Range('B636').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B638').Value = 'Cojumatlán De Régules'
$ws.Range('B703').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B728').Value = 'Coatlán Del Río'
$ws.Range('B735').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B738').Value = 'Puente De Ixtla'
$ws.Range('B743').Value = 'Tetela Del Volcán'
$ws.Range('B745').Value = 'Tlaltizapán De Zapata'
$ws.Range('B752').Value = 'Zacualpan De Amilpas'
$ws.Range('B756').Value = 'Amatlán De Cañas'
$ws.Range('B758').Value = 'Ixtlán Del Río'
$ws.Range('B765').Value = 'Santa María Del Oro'
$ws.Range('B785').Value = 'Mier Y Noriega'
$ws.Range('B790').Value = 'San Nicolás De Los Garza'
$ws.Range('B797').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B807').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B810').Value = 'Coicoyán De Las Flores'
$ws.Range('B812').Value = 'Constancia Del Rosario'
$ws.Range('B814').Value = 'El Barrio De La Soledad'
$ws.Range('B815').Value = 'Fresnillo De Trujano'
$ws.Range('B817').Value = 'Guevea De Humboldt'
$ws.Range('B818').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B819').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B820').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B821').Value = 'Huautla De Jiménez'
$ws.Range('B823').Value = 'Ixtlán De Juárez'
$ws.Range('B824').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B831').Value = 'Mariscala De Juárez'
$ws.Range('B832').Value = 'Mártires De Tacubaya'
$ws.Range('B835').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B837').Value = 'Nejapa De Madero'
$ws.Range('B838').Value = 'Oaxaca De Juárez'
$ws.Range('B839').Value = 'Ocotlán De Morelos'
$ws.Range('B840').Value = 'Pinotepa De Don Luis'
$ws.Range('B842').Value = 'Putla Villa De Guerrero'
$ws.Range('B843').Value = 'Reforma De Pineda'
$ws.Range('B852').Value = 'San Antonino El Alto'
$ws.Range('B866').Value = 'San Dionisio Del Mar'
$ws.Range('B870').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B883').Value = 'San José Del Progreso'
$ws.Range('B889').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B898').Value = 'San Juan De Los Cués'
$ws.Range('B899').Value = 'San Juan Del Estado'
$ws.Range('B931').Value = 'San Miguel Del Puerto'
$ws.Range('B933').Value = 'San Miguel El Grande'
$ws.Range('B947').Value = 'San Pablo Villa De Mitla'
$ws.Range('B949').Value = 'San Pedro El Alto'
$ws.Range('B979').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B983').Value = 'Santa Inés Del Monte'
$ws.Range('B994').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B1011').Value = 'Santiago Del Río'
$ws.Range('B1036').Value = 'Santo Domingo De Morelos'
$ws.Range('B1047').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B1048').Value = 'Tanetze De Zaragoza'
$ws.Range('B1050').Value = 'Tataltepec De Valdés'
$ws.Range('B1051').Value = 'Teotitlán De Flores Magón'
$ws.Range('B1052').Value = 'Teotitlán Del Valle'
$ws.Range('B1053').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B1054').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B1055').Value = 'Tlacolula De Matamoros'
$ws.Range('B1056').Value = 'Totontepec Villa De Morelos'
$ws.Range('B1058').Value = 'Villa De Chilapa De Díaz'
$ws.Range('B1059').Value = 'Villa De Etla'
$ws.Range('B1060').Value = 'Villa De Tututepec'
$ws.Range('B1062').Value = 'Villa Sola De Vega'
$ws.Range('B1064').Value = 'Zimatlán De Álvarez'
$ws.Range('B1081').Value = 'Ayotoxco De Guerrero'
$ws.Range('B1085').Value = 'Chalchicomula De Sesma'
$ws.Range('B1091').Value = 'Chila De La Sal'
$ws.Range('B1099').Value = 'Cuetzalan Del Progreso'
$ws.Range('B1108').Value = 'Huehuetlán El Chico'
$ws.Range('B1109').Value = 'Huehuetlán El Grande'
$ws.Range('B1115').Value = 'Izúcar De Matamoros'
$ws.Range('B1121').Value = 'Los Reyes De Juárez'
$ws.Range('B1128').Value = 'Palmar De Bravo'
$ws.Range('B1149').Value = 'Tecali De Herrera'
$ws.Range('B1156').Value = 'Tepanco De López'
$ws.Range('B1157').Value = 'Tepango De Rodríguez'
$ws.Range('B1161').Value = 'Tepexi De Rodríguez'
$ws.Range('B1162').Value = 'Tetela De Ocampo'
$ws.Range('B1166').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B1177').Value = 'Tuzamapan De Galeana'
$ws.Range('B1183').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B1194').Value = 'Amealco De Bonfil'
$ws.Range('B1196').Value = 'Cadereyta De Montes'
$ws.Range('B1202').Value = 'Jalpan De Serra'
$ws.Range('B1203').Value = 'Landa De Matamoros'
$ws.Range('B1206').Value = 'Pinal De Amoles'
$ws.Range('B1209').Value = 'San Juan Del Río'
$ws.Range('B1220').Value = 'Armadillo De Los Infante'
$ws.Range('B1221').Value = 'Axtla De Terrazas'
$ws.Range('B1227').Value = 'Ciudad Del Maíz'
$ws.Range('B1238').Value = 'Mexquitic De Carmona'
$ws.Range('B1243').Value = 'San Ciro De Acosta'
$ws.Range('B1249').Value = 'Santa María Del Río'
$ws.Range('B1251').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B1259').Value = 'Tanquián De Escobedo'
$ws.Range('B1263').Value = 'Villa De Arista'
$ws.Range('B1264').Value = 'Villa De Arriaga'
$ws.Range('B1265').Value = 'Villa De Guadalupe'
$ws.Range('B1266').Value = 'Villa De La Paz'
$ws.Range('B1267').Value = 'Villa De Ramos'
$ws.Range('B1268').Value = 'Villa De Reyes'
$ws.Range('B1299').Value = 'Nacozari De García'
$ws.Range('B1313').Value = 'Jalpa De Méndez'
$ws.Range('B1345').Value = 'Soto La Marina'
$ws.Range('B1362').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B1363').Value = 'Mazatecochco De José María Morelos'
$ws.Range('B1364').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B1366').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B1368').Value = 'San Pablo Del Monte'
$ws.Range('B1372').Value = 'Tepetitla De Lardizábal'
$ws.Range('B1373').Value = 'Tetla De La Solidaridad'
$ws.Range('B1391').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B1395').Value = 'Amatlán De Los Reyes'
$ws.Range('B1405').Value = 'Boca Del Río'
$ws.Range('B1410').Value = 'Castillo De Teayo'
$ws.Range('B1412').Value = 'Cazones De Herrera'
$ws.Range('B1429').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1430').Value = 'Cosautlán De Carvajal'
$ws.Range('B1446').Value = 'Hueyapan De Ocampo'
$ws.Range('B1447').Value = 'Ignacio De La Llave'
$ws.Range('B1450').Value = 'Ixhuatlán De Madero'
$ws.Range('B1451').Value = 'Ixhuatlán Del Café'
$ws.Range('B1460').Value = 'Juchique De Ferrer'
$ws.Range('B1464').Value = 'Las Vigas De Ramírez'
$ws.Range('B1465').Value = 'Lerdo De Tejada'
$ws.Range('B1468').Value = 'Martínez De La Torre'
$ws.Range('B1470').Value = 'Medellín De Bravo'
$ws.Range('B1474').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B1483').Value = 'Ozuluama De Mascareñas'
$ws.Range('B1487').Value = 'Paso De Ovejas'
$ws.Range('B1488').Value = 'Paso Del Macho'
$ws.Range('B1492').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1500').Value = 'Sayula De Alemán'
$ws.Range('B1504').Value = 'Soledad De Doblado'
$ws.Range('B1510').Value = 'Tatahuicapan De Juárez'
$ws.Range('B1543').Value = 'Vega De Alatorre'
$ws.Range('B1552').Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range('B1553').Value = 'Zozocolco De Hidalgo'
$ws.Range('B1567').Value = 'Concepción Del Oro'
$ws.Range('B1576').Value = 'Jiménez Del Teul'
$ws.Range('B1585').Value = 'Moyahua De Estrada'
$ws.Range('B1586').Value = 'Nochistlán De Mejía'
$ws.Range('B1587').Value = 'Noria De Ángeles'
$ws.Range('B1597').Value = 'Teúl De González Ortega'
$ws.Range('B1598').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1601').Value = 'Villa De Cos'

# --- Fix 1-ULP floating point differences in percentage column ---
$ws.Range('D153').Value = 0.0009713264433910948
$ws.Range('D244').Value = 0.0009324733856554512
$ws.Range('D292').Value = 0.0009324733856554512
$ws.Range('D360').Value = 0.0009324733856554512
$ws.Range('D371').Value = 0.09142124485196984
$ws.Range('D447').Value = 0.0009324733856554512
$ws.Range('D473').Value = 0.0009324733856554512
$ws.Range('D489').Value = 0.0009713264433910948
$ws.Range('D669').Value = 0.009674411376175304
$ws.Range('D680').Value = 0.0009713264433910948
$ws.Range('D759').Value = 0.0009713264433910948
$ws.Range('D768').Value = 0.0009713264433910948
$ws.Range('D847').Value = 0.0009713264433910948
$ws.Range('D956').Value = 0.0009324733856554512
$ws.Range('D1078').Value = 0.0009713264433910948
$ws.Range('D1216').Value = 0.0009324733856554512
$ws.Range('D1218').Value = 0.0009324733856554512
$ws.Range('D1570').Value = 0.0009324733856554512

# --- Remove trailing footer/metadata rows (source note, sample size, etc.) ---
$ws.Range("A1610:A1614").EntireRow.Delete()

